# Apply "Errors fixed and updated" changes to ProductRules workbook.
$wb = $excel.ActiveWorkbook

$response = $wb.Worksheets.Item("Response")

# --- Update data values on the "Response" sheet ---
# Child_minCoverAllowed for TC_004_1960E row: 5000 -> 10000
$response.Range("X5").Value = 10000

# TC_005_1960C row: productName Elite V6.0 -> Core V6.0, Mainlife_maxAge 59 -> 75
$response.Range("D6").Value = "Core V6.0"
$response.Range("M6").Value = 75

# TC_006_1960L row: productName Elite V6.0 -> Lite V6.0
$response.Range("D7").Value = "Lite V6.0"

# --- Update the active sheet / selection so "Response" becomes the visible tab ---
$response.Activate()
$response.Range("D7").Select()
